$d = $word.ActiveDocument
$d.Content.Find.Execute("CONTRAT DE CAUTIONNEMENT", $true, $false, $false, $false, $false, $true, 1, $false, "CONTRAT DE CAUTIONNEMENT", 2)
